$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 12
$ws.Range("B3").Value = 22
$ws.Range("C3").Value = 22
$ws.Range("D3").Value = 22
$ws.Range("E3").Value = 22
